# "Completed basic structure UwU"
# Update the candidate Status column and fix a mismatched "Current Position"
# value, then leave the selection where the author last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status column (H): reword the placeholder statuses ---
# Every cell that used to read "YES" becomes "To be checked"
$ws.Range("H2").Value = "To be checked"
$ws.Range("H3").Value = "To be checked"
$ws.Range("H4").Value = "To be checked"
$ws.Range("H6").Value = "To be checked"
# The single "NO" becomes "Good"
$ws.Range("H5").Value = "Good"

# --- Current Position column (G): row 6 (Amanda Dieudonné) had the wrong
#     current position copied over; correct it to match the submitted
#     position ("Operation Executive Assistant") ---
$ws.Range("G6").Value = "Operation Executive Assistant"

# --- Leave the selection on G14, matching where the author ended up ---
$ws.Range("G14").Select() | Out-Null
